# "search priority by table"
# - Existing "Лист1" schedule sheet: bump the teachers' search-priority
#   numbers in column F (and add one more subject for "степан").
# - New sheet "Лист2" inserted before "Лист1": a small scratch/search table.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update the existing schedule sheet (stays named "Лист1") ---
$ws1.Range("C2").Value = "информатика, математика,физика"
$ws1.Range("F2").Value = 2
$ws1.Range("F3").Value = 2
$ws1.Range("F4").Value = 3

# Keep the rest (F5 priority, rows 7-12, etc.) untouched.

# Set the selection on "Лист1" now, while it is still the active sheet --
# selections only persist to a sheet's saved view state while that sheet
# is the active one at the time .Select() runs.
$ws1.Range("F6").Select()

# --- Insert the new "Лист2" sheet before "Лист1" ---
$ws2 = $wb.Worksheets.Add()
$ws2.Name = "Лист2"

$ws2.Range("A1").Value = ",vp[av"
$ws2.Range("B1").Value = "avdf"
$ws2.Range("C1").Value = "dfvfdv"
$ws2.Range("C2").Value = "avddfv"

# "Лист2" becomes the active sheet on Add(); set its selection too.
$ws2.Range("C2").Select()
